$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing row (66) down to the new row (67)
# so the new row picks up the same "Medium" fill on column B and the same
# Hyperlink style on column E without fabricating brand new styles.
$ws.Range("A66:E66").Copy() | Out-Null
$ws.Range("A67:E67").PasteSpecial(-4122) | Out-Null

# New question row content
$ws.Range("A67").Value = "2300. Successful Pairs of Spells and Potions"
$ws.Range("B67").Value = "Medium"
$ws.Range("C67").Value = "Binary Search"
$ws.Range("D67").Value = "First sort the pairs array. For loop over all spells and do a binary search to find the first mid point which results in a success. Because we previously sorted potions, all indices after that are valid, so add that count to the pairs[i]."
$ws.Range("E67").Value = "https://leetcode.com/problems/successful-pairs-of-spells-and-potions/solutions/3367914/easy-solutions-in-java-python-and-c-look-at-once-with-exaplanation/?envType=study-plan-v2&envId=leetcode-75 "

# Wire up the hyperlink for the new Link cell, then restore the Hyperlink
# cell style (Add() re-applies font formatting as a fresh style, so reset
# it back to the shared "Hyperlink" style used throughout the sheet).
$ws.Hyperlinks.Add($ws.Range("E67"), "https://leetcode.com/problems/successful-pairs-of-spells-and-potions/solutions/3367914/easy-solutions-in-java-python-and-c-look-at-once-with-exaplanation/?envType=study-plan-v2&envId=leetcode-75 ") | Out-Null
$ws.Range("E67").Style = "Hyperlink"

# Update selection/active cell to match post-edit state
$ws.Range("E73").Select() | Out-Null
